$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine rows 2-7 into a single Python-tuple-like string in A2
$combined = "('תהילה', ['{3}{W}{W}', 'יצור — התגלמות', 'יכולת תﬠופה', '{W}{2}: יצורים שנמצאים בשליטתך נהנים מהגנתו של הצבﬠ שבחרת, ﬠד לסיום החור. השתמש/י ביכולת זו רק במקרה, שתהילה נמצאת בבית- הקברות שלך.', '3/3'])"

$ws.Range("A2").Value = $combined

# Remove now-obsolete rows 3-7
$ws.Range("A3:A7").EntireRow.Delete()
